$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text fields (row header labels stay the same, only values in column B change)
$ws.Range("B2").Value = "Petros Petropoulos"

# B3's new value "69" looks numeric, so force it to stay text using the
# leading-apostrophe convention (same as typing '69 into Excel directly).
$ws.Range("B3").Formula = "'69"

$ws.Range("B4").Value = "R_SWING (2)"
$ws.Range("B5").Value = "jlj"
$ws.Range("B6").Value = "jlkj"
$ws.Range("B7").Value = "jlkj"

# Numeric fields
$ws.Range("B8").Value = 5
$ws.Range("B9").Value = 5
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 50

$ws.Range("B13").Value = 45
$ws.Range("B14").Value = 565

# B15's new value "dfdf" is plain text
$ws.Range("B15").Value = "dfdf"

$ws.Range("B16").Value = 5
$ws.Range("B17").Value = 6
$ws.Range("B18").Value = 8
$ws.Range("B19").Value = 4
$ws.Range("B20").Value = 5
$ws.Range("B21").Value = 68
$ws.Range("B22").Value = 979
$ws.Range("B23").Value = 7
$ws.Range("B24").Value = 56
$ws.Range("B25").Value = 76
$ws.Range("B26").Value = 5
$ws.Range("B27").Value = 75
$ws.Range("B28").Value = 675
$ws.Range("B29").Value = 65
$ws.Range("B30").Value = 765

$ws.Range("B32").Value = 654
$ws.Range("B33").Value = 86
$ws.Range("B34").Value = 78
$ws.Range("B35").Value = 6785
$ws.Range("B36").Value = 47
$ws.Range("B37").Value = 687
$ws.Range("B38").Value = 58
$ws.Range("B39").Value = 67
$ws.Range("B40").Value = 57
$ws.Range("B41").Value = 857
$ws.Range("B42").Value = 686
$ws.Range("B43").Value = 7
$ws.Range("B44").Value = 657
$ws.Range("B45").Value = 65
$ws.Range("B46").Value = 7
$ws.Range("B47").Value = 568
